$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 192,2
$arr[0,0] = 45984; $arr[0,1] = 454
$arr[1,0] = 45984.01041666666; $arr[1,1] = 449
$arr[2,0] = 45984.02083333334; $arr[2,1] = 455
$arr[3,0] = 45984.03125; $arr[3,1] = 451
$arr[4,0] = 45984.04166666666; $arr[4,1] = 447
$arr[5,0] = 45984.05208333334; $arr[5,1] = 445
$arr[6,0] = 45984.0625; $arr[6,1] = 0
$arr[7,0] = 45984.07291666666; $arr[7,1] = 0
$arr[8,0] = 45984.08333333334; $arr[8,1] = 450
$arr[9,0] = 45984.09375; $arr[9,1] = 0
$arr[10,0] = 45984.10416666666; $arr[10,1] = 0
$arr[11,0] = 45984.11458333334; $arr[11,1] = 0
$arr[12,0] = 45984.125; $arr[12,1] = 465
$arr[13,0] = 45984.13541666666; $arr[13,1] = 464
$arr[14,0] = 45984.14583333334; $arr[14,1] = 0
$arr[15,0] = 45984.15625; $arr[15,1] = 466
$arr[16,0] = 45984.16666666666; $arr[16,1] = 473
$arr[17,0] = 45984.17708333334; $arr[17,1] = 472
$arr[18,0] = 45984.1875; $arr[18,1] = 473
$arr[19,0] = 45984.19791666666; $arr[19,1] = 475
$arr[20,0] = 45984.20833333334; $arr[20,1] = 507
$arr[21,0] = 45984.21875; $arr[21,1] = 508
$arr[22,0] = 45984.22916666666; $arr[22,1] = 509
$arr[23,0] = 45984.23958333334; $arr[23,1] = 535
$arr[24,0] = 45984.25; $arr[24,1] = 591
$arr[25,0] = 45984.26041666666; $arr[25,1] = 573
$arr[26,0] = 45984.27083333334; $arr[26,1] = 562
$arr[27,0] = 45984.28125; $arr[27,1] = 560
$arr[28,0] = 45984.29166666666; $arr[28,1] = 603
$arr[29,0] = 45984.30208333334; $arr[29,1] = 610
$arr[30,0] = 45984.3125; $arr[30,1] = 626
$arr[31,0] = 45984.32291666666; $arr[31,1] = 629
$arr[32,0] = 45984.33333333334; $arr[32,1] = 634
$arr[33,0] = 45984.34375; $arr[33,1] = 639
$arr[34,0] = 45984.35416666666; $arr[34,1] = 640
$arr[35,0] = 45984.36458333334; $arr[35,1] = 639
$arr[36,0] = 45984.375; $arr[36,1] = 611
$arr[37,0] = 45984.38541666666; $arr[37,1] = 606
$arr[38,0] = 45984.39583333334; $arr[38,1] = 607
$arr[39,0] = 45984.40625; $arr[39,1] = 606
$arr[40,0] = 45984.41666666666; $arr[40,1] = 628
$arr[41,0] = 45984.42708333334; $arr[41,1] = 789
$arr[42,0] = 45984.4375; $arr[42,1] = 798
$arr[43,0] = 45984.44791666666; $arr[43,1] = 789
$arr[44,0] = 45984.45833333334; $arr[44,1] = 556
$arr[45,0] = 45984.46875; $arr[45,1] = 552
$arr[46,0] = 45984.47916666666; $arr[46,1] = 553
$arr[47,0] = 45984.48958333334; $arr[47,1] = 550
$arr[48,0] = 45984.5; $arr[48,1] = 671
$arr[49,0] = 45984.51041666666; $arr[49,1] = 736
$arr[50,0] = 45984.52083333334; $arr[50,1] = 741
$arr[51,0] = 45984.53125; $arr[51,1] = 746
$arr[52,0] = 45984.54166666666; $arr[52,1] = 716
$arr[53,0] = 45984.55208333334; $arr[53,1] = 724
$arr[54,0] = 45984.5625; $arr[54,1] = 0
$arr[55,0] = 45984.57291666666; $arr[55,1] = 723
$arr[56,0] = 45984.58333333334; $arr[56,1] = 726
$arr[57,0] = 45984.59375; $arr[57,1] = 640
$arr[58,0] = 45984.60416666666; $arr[58,1] = 632
$arr[59,0] = 45984.61458333334; $arr[59,1] = 637
$arr[60,0] = 45984.625; $arr[60,1] = 730
$arr[61,0] = 45984.63541666666; $arr[61,1] = 732
$arr[62,0] = 45984.64583333334; $arr[62,1] = 734
$arr[63,0] = 45984.65625; $arr[63,1] = 736
$arr[64,0] = 45984.66666666666; $arr[64,1] = 784
$arr[65,0] = 45984.67708333334; $arr[65,1] = 789
$arr[66,0] = 45984.6875; $arr[66,1] = 793
$arr[67,0] = 45984.69791666666; $arr[67,1] = 825
$arr[68,0] = 45984.70833333334; $arr[68,1] = 1063
$arr[69,0] = 45984.71875; $arr[69,1] = 1064
$arr[70,0] = 45984.72916666666; $arr[70,1] = 1070
$arr[71,0] = 45984.73958333334; $arr[71,1] = 1150
$arr[72,0] = 45984.75; $arr[72,1] = 1180
$arr[73,0] = 45984.76041666666; $arr[73,1] = 1155
$arr[74,0] = 45984.77083333334; $arr[74,1] = 1111
$arr[75,0] = 45984.78125; $arr[75,1] = 1110
$arr[76,0] = 45984.79166666666; $arr[76,1] = 1057
$arr[77,0] = 45984.80208333334; $arr[77,1] = 1025
$arr[78,0] = 45984.8125; $arr[78,1] = 1022
$arr[79,0] = 45984.82291666666; $arr[79,1] = 1024
$arr[80,0] = 45984.83333333334; $arr[80,1] = 930
$arr[81,0] = 45984.84375; $arr[81,1] = 846
$arr[82,0] = 45984.85416666666; $arr[82,1] = 868
$arr[83,0] = 45984.86458333334; $arr[83,1] = 977
$arr[84,0] = 45984.875; $arr[84,1] = 923
$arr[85,0] = 45984.88541666666; $arr[85,1] = 913
$arr[86,0] = 45984.89583333334; $arr[86,1] = 938
$arr[87,0] = 45984.90625; $arr[87,1] = 929
$arr[88,0] = 45984.91666666666; $arr[88,1] = 573
$arr[89,0] = 45984.92708333334; $arr[89,1] = 563
$arr[90,0] = 45984.9375; $arr[90,1] = 564
$arr[91,0] = 45984.94791666666; $arr[91,1] = 558
$arr[92,0] = 45984.95833333334; $arr[92,1] = 478
$arr[93,0] = 45984.96875; $arr[93,1] = 474
$arr[94,0] = 45984.97916666666; $arr[94,1] = 477
$arr[95,0] = 45984.98958333334; $arr[95,1] = 476
$arr[96,0] = 45985; $arr[96,1] = 462
$arr[97,0] = 45985.01041666666; $arr[97,1] = 475
$arr[98,0] = 45985.02083333334; $arr[98,1] = 476
$arr[99,0] = 45985.03125; $arr[99,1] = 0
$arr[100,0] = 45985.04166666666; $arr[100,1] = 477
$arr[101,0] = 45985.05208333334; $arr[101,1] = 473
$arr[102,0] = 45985.0625; $arr[102,1] = 472
$arr[103,0] = 45985.07291666666; $arr[103,1] = 473
$arr[104,0] = 45985.08333333334; $arr[104,1] = 468
$arr[105,0] = 45985.09375; $arr[105,1] = 0
$arr[106,0] = 45985.10416666666; $arr[106,1] = 470
$arr[107,0] = 45985.11458333334; $arr[107,1] = 468
$arr[108,0] = 45985.125; $arr[108,1] = 471
$arr[109,0] = 45985.13541666666; $arr[109,1] = 470
$arr[110,0] = 45985.14583333334; $arr[110,1] = 472
$arr[111,0] = 45985.15625; $arr[111,1] = 478
$arr[112,0] = 45985.16666666666; $arr[112,1] = 533
$arr[113,0] = 45985.17708333334; $arr[113,1] = 532
$arr[114,0] = 45985.1875; $arr[114,1] = 530
$arr[115,0] = 45985.19791666666; $arr[115,1] = 536
$arr[116,0] = 45985.20833333334; $arr[116,1] = 548
$arr[117,0] = 45985.21875; $arr[117,1] = 545
$arr[118,0] = 45985.22916666666; $arr[118,1] = 548
$arr[119,0] = 45985.23958333334; $arr[119,1] = 572
$arr[120,0] = 45985.25; $arr[120,1] = 950
$arr[121,0] = 45985.26041666666; $arr[121,1] = 992
$arr[122,0] = 45985.27083333334; $arr[122,1] = 1004
$arr[123,0] = 45985.28125; $arr[123,1] = 0
$arr[124,0] = 45985.29166666666; $arr[124,1] = 1040
$arr[125,0] = 45985.30208333334; $arr[125,1] = 1041
$arr[126,0] = 45985.3125; $arr[126,1] = 1043
$arr[127,0] = 45985.32291666666; $arr[127,1] = 1052
$arr[128,0] = 45985.33333333334; $arr[128,1] = 1102
$arr[129,0] = 45985.34375; $arr[129,1] = 1113
$arr[130,0] = 45985.35416666666; $arr[130,1] = 1103
$arr[131,0] = 45985.36458333334; $arr[131,1] = 1199
$arr[132,0] = 45985.375; $arr[132,1] = 1034
$arr[133,0] = 45985.38541666666; $arr[133,1] = 1100
$arr[134,0] = 45985.39583333334; $arr[134,1] = 1163
$arr[135,0] = 45985.40625; $arr[135,1] = 1168
$arr[136,0] = 45985.41666666666; $arr[136,1] = 1134
$arr[137,0] = 45985.42708333334; $arr[137,1] = 1120
$arr[138,0] = 45985.4375; $arr[138,1] = 979
$arr[139,0] = 45985.44791666666; $arr[139,1] = 983
$arr[140,0] = 45985.45833333334; $arr[140,1] = 703
$arr[141,0] = 45985.46875; $arr[141,1] = 698
$arr[142,0] = 45985.47916666666; $arr[142,1] = 699
$arr[143,0] = 45985.48958333334; $arr[143,1] = 696
$arr[144,0] = 45985.5; $arr[144,1] = 683
$arr[145,0] = 45985.51041666666; $arr[145,1] = 682
$arr[146,0] = 45985.52083333334; $arr[146,1] = 684
$arr[147,0] = 45985.53125; $arr[147,1] = 0
$arr[148,0] = 45985.54166666666; $arr[148,1] = 0
$arr[149,0] = 45985.55208333334; $arr[149,1] = 0
$arr[150,0] = 45985.5625; $arr[150,1] = 0
$arr[151,0] = 45985.57291666666; $arr[151,1] = 0
$arr[152,0] = 45985.58333333334; $arr[152,1] = 0
$arr[153,0] = 45985.59375; $arr[153,1] = 0
$arr[154,0] = 45985.60416666666; $arr[154,1] = 0
$arr[155,0] = 45985.61458333334; $arr[155,1] = 0
$arr[156,0] = 45985.625; $arr[156,1] = 0
$arr[157,0] = 45985.63541666666; $arr[157,1] = 0
$arr[158,0] = 45985.64583333334; $arr[158,1] = 0
$arr[159,0] = 45985.65625; $arr[159,1] = 0
$arr[160,0] = 45985.66666666666; $arr[160,1] = 0
$arr[161,0] = 45985.67708333334; $arr[161,1] = 0
$arr[162,0] = 45985.6875; $arr[162,1] = 0
$arr[163,0] = 45985.69791666666; $arr[163,1] = 0
$arr[164,0] = 45985.70833333334; $arr[164,1] = 0
$arr[165,0] = 45985.71875; $arr[165,1] = 0
$arr[166,0] = 45985.72916666666; $arr[166,1] = 0
$arr[167,0] = 45985.73958333334; $arr[167,1] = 0
$arr[168,0] = 45985.75; $arr[168,1] = 0
$arr[169,0] = 45985.76041666666; $arr[169,1] = 0
$arr[170,0] = 45985.77083333334; $arr[170,1] = 0
$arr[171,0] = 45985.78125; $arr[171,1] = 0
$arr[172,0] = 45985.79166666666; $arr[172,1] = 0
$arr[173,0] = 45985.80208333334; $arr[173,1] = 0
$arr[174,0] = 45985.8125; $arr[174,1] = 0
$arr[175,0] = 45985.82291666666; $arr[175,1] = 0
$arr[176,0] = 45985.83333333334; $arr[176,1] = 0
$arr[177,0] = 45985.84375; $arr[177,1] = 0
$arr[178,0] = 45985.85416666666; $arr[178,1] = 0
$arr[179,0] = 45985.86458333334; $arr[179,1] = 0
$arr[180,0] = 45985.875; $arr[180,1] = 0
$arr[181,0] = 45985.88541666666; $arr[181,1] = 0
$arr[182,0] = 45985.89583333334; $arr[182,1] = 0
$arr[183,0] = 45985.90625; $arr[183,1] = 0
$arr[184,0] = 45985.91666666666; $arr[184,1] = 0
$arr[185,0] = 45985.92708333334; $arr[185,1] = 0
$arr[186,0] = 45985.9375; $arr[186,1] = 0
$arr[187,0] = 45985.94791666666; $arr[187,1] = 0
$arr[188,0] = 45985.95833333334; $arr[188,1] = 0
$arr[189,0] = 45985.96875; $arr[189,1] = 0
$arr[190,0] = 45985.97916666666; $arr[190,1] = 0
$arr[191,0] = 45985.98958333334; $arr[191,1] = 0

$ws.Range("A2:B193").Value = $arr
